$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: B9:K9 - fill in availability values
$ws.Range("B9").Value = ":("
$ws.Range("C9").Value = ":D"
$ws.Range("D9").Value = ":("
$ws.Range("E9").Value = ":D"
$ws.Range("F9").Value = ":D"
$ws.Range("G9").Value = ":D"
$ws.Range("H9").Value = ":("
$ws.Range("I9").Value = ":D"
$ws.Range("J9").Value = ":D - :("
$ws.Range("K9").Value = ":D"

# Row 17: B17:K17 - fill in availability values
$ws.Range("B17").Value = ":D"
$ws.Range("C17").Value = ":D"
$ws.Range("D17").Value = ":("
$ws.Range("E17").Value = ":D"
$ws.Range("F17").Value = ":("
$ws.Range("G17").Value = ":D"
$ws.Range("H17").Value = ":("
$ws.Range("I17").Value = ":D"
$ws.Range("J17").Value = ":("
$ws.Range("K17").Value = ":D"

# Update the selection to match the target (J17 was last edited cell)
$ws.Range("J17").Select()
